$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated roster data (rows 2-18), columns: A=Player, B=Position, C=Team
$data = @(
    @("Lonzo Ball", "PG", "Chicago Bulls"),
    @("Donovan Mitchell", "PG,SG", "Cleveland Cavaliers"),
    @("Malik Beasley", "SG,SF", "Detroit Pistons"),
    @("Josh Hart", "SG,SF,PF", "New York Knicks"),
    @("Andrew Wiggins", "SF,PF", "Golden State Warriors"),
    @("De'Andre Hunter", "SF,PF", "Atlanta Hawks"),
    @("Dyson Daniels", "PG,SG,SF", "Atlanta Hawks"),
    @("Santi Aldama", "PF,C", "Memphis Grizzlies"),
    @("Alperen Sengün", "C", "Houston Rockets"),
    @("Kristaps Porzingis", "PF,C", "Boston Celtics"),
    @("Brice Sensabaugh", "SF,PF", "Utah Jazz"),
    @("Domantas Sabonis", "C", "Sacramento Kings"),
    @("Victor Wembanyama", "C", "San Antonio Spurs"),
    @("Kelly Oubre Jr.", "SG,SF", "Philadelphia 76ers"),
    @("Michael Porter Jr.", "SF,PF", "Denver Nuggets"),
    @("Cam Thomas", "SG,SF", "Brooklyn Nets"),
    @("Donte DiVincenzo", "PG,SG,SF", "Minnesota Timberwolves")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
